$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update an existing tag: row 61, column C (interlocking) "p" -> "d" ---
$ws.Range("C61").Value = "d"

# --- Append new keyword rows 217-221 in columns A (keyword) and C (tag) ---
$ws.Range("A217").Value = "stick on"
$ws.Range("C217").Value = "p"

$ws.Range("A218").Value = "waterproofing"
$ws.Range("C218").Value = "s"

$ws.Range("A219").Value = "terracotta"
$ws.Range("C219").Value = "d"

$ws.Range("A220").Value = "deck"
$ws.Range("C220").Value = "a"

$ws.Range("A221").Value = "foam"
$ws.Range("C221").Value = "m"

# --- Match the author's final on-screen selection / scroll position ---
$ws.Range("B221").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 202
$win.ScrollColumn = 1

$wb.Save()
